# Daily roll-forward update for the "剩余" (days remaining) tracker.
# Rule (one day has elapsed since the workbook was last generated):
#   - For every data row, the "剩余" (remaining days, column E) counts down by 1.
#   - If E would hit 0 (i.e. E was 1), the billing cycle renews instead:
#       E resets to the row's total-day count (column D), and the
#       "开始时间" start date (column F, an yyyyMMdd integer) advances by
#       D calendar days.
#   - Rows whose start date isn't a well-formed yyyyMMdd value (data glitch)
#     are left untouched, since the date math can't be performed on them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    if ($d -eq $null -or $e -eq $null) { continue }

    $fstr = [string]$f
    if ($fstr.Length -ne 8) {
        # malformed start date (e.g. "202510929") - skip, matches source behavior
        continue
    }

    if ($e -eq 1) {
        $newE = $d
        $dt = [datetime]::ParseExact($fstr, "yyyyMMdd", $null)
        $dt2 = $dt.AddDays($d)
        $newF = [int]$dt2.ToString("yyyyMMdd")
    } else {
        $newE = $e - 1
        $newF = $f
    }

    $ws.Cells.Item($r, 5).Value2 = $newE
    $ws.Cells.Item($r, 6).Value2 = $newF
}
